$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "amount"
$ws.Range("C2").Value = 1500
